$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:J12 (rows 2-12, columns B-J)
$data = @(
    @(13, 1400, 1504, 17, 4, 23, 72, 1501.5, 0.166500166500172),
    @(1, 136, 145, 7, 1, 2, 0, 148, -2.027027027027029),
    @(3, 327, 351, 20, 2, 4, 0, 120, 192.5),
    @(1, 789, 829, 5, 3, 18, 33, 926, -10.47516198704104),
    @(2, 627, 707, 67, 10, 9, 1, 642, 10.12461059190031),
    @(1, 186, 190, 3, 0, 0, 0, 196, -3.061224489795922),
    @(0, 89, 113, 21, 3, 0, 0, 119, -5.042016806722693),
    @(0, 364, 364, 0, 0, 9, 0, 788, -53.80710659898477),
    @(2, 78, 80, 0, 0, 2, 0, 95, -15.78947368421053),
    @(0, 3, 3, 0, 0, 0, 0, 4, -25),
    @(1, 47, 48, 0, 0, 4, 0, 85, -43.52941176470588)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
